# Auto-generated Excel COM-interop edit script.
# Applies the numeric cell updates described by the commit diff to the
# underlying FFXIV leve-profit workbook. Touches the ALC / ARM / BSM / CRP /
# GSM / LTW / WVR sheets (CUL has no changes in the diff).
#
# All touched cells are plain numeric literals (no formulas anywhere in the
# workbook), so each one is just a direct Range.Value assignment. A few
# cells are deleted outright by the diff (the M/N "profit" column was blank
# before the edit touched neighbouring rows) -- those use ClearContents()
# so the saved XML drops the <c> element entirely, matching how every other
# blank cell in this workbook is represented (no <c> at all, not a 0).

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 4000
$ws.Range("J13").Value = 4000
$ws.Range("L13").Value = 4000
$ws.Range("N13").Value = -4338
$ws.Range("H70").Value = 2333.3333
$ws.Range("I70").Value = 2320
$ws.Range("J70").Value = 2350
$ws.Range("K70").Value = 6960
$ws.Range("L70").Value = 7050
$ws.Range("M70").Value = -6690
$ws.Range("N70").Value = -7590
$ws.Range("H73").Value = 2333.3333
$ws.Range("I73").Value = 2320
$ws.Range("J73").Value = 2350
$ws.Range("K73").Value = 6960
$ws.Range("L73").Value = 7050
$ws.Range("M73").Value = -6024
$ws.Range("N73").Value = -8922
$ws.Range("H80").Value = 1336
$ws.Range("I80").Value = 651
$ws.Range("J80").Value = 1610
$ws.Range("K80").Value = 1953
$ws.Range("L80").Value = 4830
$ws.Range("M80").Value = -955
$ws.Range("N80").Value = -6826
$ws.Range("H83").Value = 1336
$ws.Range("I83").Value = 651
$ws.Range("J83").Value = 1610
$ws.Range("K83").Value = 5859
$ws.Range("L83").Value = 14490
$ws.Range("M83").Value = -867
$ws.Range("N83").Value = -24474
$ws.Range("H100").Value = 62501484
$ws.Range("I100").Value = 1720
$ws.Range("J100").Value = 166667760
$ws.Range("K100").Value = 1720
$ws.Range("L100").Value = 166667760
$ws.Range("M100").Value = -1179
$ws.Range("N100").Value = -166668842
$ws.Range("H101").Value = 1735.6666
$ws.Range("I101").Value = 1484.3636
$ws.Range("J101").Value = 4500
$ws.Range("K101").Value = 4453.0908
$ws.Range("L101").Value = 13500
$ws.Range("M101").Value = -2831.0908
$ws.Range("N101").Value = -16744
$ws.Range("H137").Value = 379095.66
$ws.Range("I137").Value = 734200.9
$ws.Range("J137").Value = 6235.15
$ws.Range("K137").Value = 2202602.7
$ws.Range("L137").Value = 18705.45
$ws.Range("M137").Value = -2200052.7
$ws.Range("N137").Value = -23805.45

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 501.16666
$ws.Range("I26").Value = 501.16666
$ws.Range("K26").Value = 501.16666
$ws.Range("M26").Value = -171.16666
$ws.Range("H27").Value = 16899.666
$ws.Range("J27").Value = 16899.666
$ws.Range("L27").Value = 16899.666
$ws.Range("N27").Value = -17267.666
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()  # was -1197, now blank
$ws.Range("H80").Value = 31452.5
$ws.Range("J80").Value = 38903.332
$ws.Range("L80").Value = 38903.332
$ws.Range("N80").Value = -40899.332
$ws.Range("H83").Value = 31452.5
$ws.Range("J83").Value = 38903.332
$ws.Range("L83").Value = 116709.996
$ws.Range("N83").Value = -126693.996

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H49").Value = 43000
$ws.Range("J49").Value = 43000
$ws.Range("L49").Value = 43000
$ws.Range("N49").Value = -43478
$ws.Range("H130").Value = 51811.43
$ws.Range("J130").Value = 51811.43
$ws.Range("L130").Value = 51811.43
$ws.Range("N130").Value = -61851.43

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2971.0337
$ws.Range("I31").Value = 1436.8718
$ws.Range("J31").Value = 4167.68
$ws.Range("K31").Value = 1436.8718
$ws.Range("L31").Value = 4167.68
$ws.Range("M31").Value = -1141.8718
$ws.Range("N31").Value = -4757.68
$ws.Range("H34").Value = 2971.0337
$ws.Range("I34").Value = 1436.8718
$ws.Range("J34").Value = 4167.68
$ws.Range("K34").Value = 1436.8718
$ws.Range("L34").Value = 4167.68
$ws.Range("M34").Value = -1234.8718
$ws.Range("N34").Value = -4571.68
$ws.Range("H81").Value = 34000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 34000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 34000
$ws.Range("M81").ClearContents()  # was -29002, now blank
$ws.Range("N81").Value = -35996
$ws.Range("H84").Value = 34000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 34000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 102000
$ws.Range("M84").ClearContents()  # was -85008, now blank
$ws.Range("N84").Value = -111984
$ws.Range("H131").Value = 40856.8
$ws.Range("J131").Value = 40856.8
$ws.Range("L131").Value = 40856.8
$ws.Range("N131").Value = -50936.8

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I113").Value = 1985.4546
$ws.Range("J113").Value = 2122.5
$ws.Range("K113").Value = 1985.4546
$ws.Range("L113").Value = 2122.5
$ws.Range("M113").Value = 184.5454
$ws.Range("N113").Value = -6462.5
$ws.Range("H126").Value = 6955.1816
$ws.Range("I126").Value = 3685.7144
$ws.Range("J126").Value = 8480.933999999999
$ws.Range("K126").Value = 11057.1432
$ws.Range("L126").Value = 25442.802
$ws.Range("M126").Value = -8587.143199999999
$ws.Range("N126").Value = -30382.802

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 618.3684
$ws.Range("I22").Value = 341
$ws.Range("J22").Value = 999.75
$ws.Range("K22").Value = 341
$ws.Range("L22").Value = 999.75
$ws.Range("M22").Value = -46
$ws.Range("N22").Value = -1589.75
$ws.Range("H27").Value = 618.3684
$ws.Range("I27").Value = 341
$ws.Range("J27").Value = 999.75
$ws.Range("K27").Value = 341
$ws.Range("L27").Value = 999.75
$ws.Range("M27").Value = -234
$ws.Range("N27").Value = -1213.75
$ws.Range("H82").Value = 2298.5
$ws.Range("I82").Value = 1700.3334
$ws.Range("J82").Value = 2896.6667
$ws.Range("K82").Value = 1700.3334
$ws.Range("L82").Value = 2896.6667
$ws.Range("M82").Value = -1339.3334
$ws.Range("N82").Value = -3618.6667
$ws.Range("H85").Value = 2298.5
$ws.Range("I85").Value = 1700.3334
$ws.Range("J85").Value = 2896.6667
$ws.Range("K85").Value = 1700.3334
$ws.Range("L85").Value = 2896.6667
$ws.Range("M85").Value = -452.3334
$ws.Range("N85").Value = -5392.6667
$ws.Range("H93").Value = 2181.111
$ws.Range("I93").Value = 1953.75
$ws.Range("K93").Value = 1953.75
$ws.Range("M93").Value = -705.75

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 40350.445
$ws.Range("J82").Value = 40350.445
$ws.Range("L82").Value = 40350.445
$ws.Range("N82").Value = -41116.445
$ws.Range("H85").Value = 40350.445
$ws.Range("J85").Value = 40350.445
$ws.Range("L85").Value = 40350.445
$ws.Range("N85").Value = -43002.445
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()  # was -27214, now blank
$ws.Range("H132").Value = 1913.8358
$ws.Range("I132").Value = 1724.38
$ws.Range("J132").Value = 2471.0588
$ws.Range("K132").Value = 5173.14
$ws.Range("L132").Value = 7413.176399999999
$ws.Range("M132").Value = -2643.14
$ws.Range("N132").Value = -12473.1764

Write-Output "Applied 177 cell updates across 7 sheets"
